$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.006090927296766
$ws.Range("D2").Value = 1.008736917786339
$ws.Range("E2").Value = 1.008722147100834
$ws.Range("F2").Value = 1.004255285924022
$ws.Range("I2").Value = 1.023594999628091
$ws.Range("J2").Value = 1.011370418211365
$ws.Range("K2").Value = 1.011612260904118
$ws.Range("L2").Value = 1.011597535032606
$ws.Range("M2").Value = 1.00714429220932
$ws.Range("N2").Value = 1.007347764093004
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.008747479273435
$ws.Range("D3").Value = 1.011303855525857
$ws.Range("E3").Value = 1.011042418778869
$ws.Range("F3").Value = 1.007591720365494
$ws.Range("I3").Value = 1.023504579208683
$ws.Range("J3").Value = 1.013648668426333
$ws.Range("K3").Value = 1.01397794869283
$ws.Range("L3").Value = 1.013717244929004
$ws.Range("M3").Value = 1.010276259819717
$ws.Range("N3").Value = 1.008163843152681
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.010454205672261
$ws.Range("D4").Value = 1.012953150827621
$ws.Range("E4").Value = 1.012532685184043
$ws.Range("F4").Value = 1.009736437635741
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.015110944212857
$ws.Range("K4").Value = 1.015496830735572
$ws.Range("L4").Value = 1.015077481814618
$ws.Range("M4").Value = 1.012288686452855
$ws.Range("N4").Value = 1.008686106013879
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.011168864086611
$ws.Range("D5").Value = 1.01364379410684
$ws.Range("E5").Value = 1.013156603907741
$ws.Range("F5").Value = 1.010634790494188
$ws.Range("I5").Value = 1.023416932628352
$ws.Range("J5").Value = 1.015722906349846
$ws.Range("K5").Value = 1.016132598703526
$ws.Range("L5").Value = 1.015646672322158
$ws.Range("M5").Value = 1.013131421759706
$ws.Range("N5").Value = 1.00890430323948
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.01128869361613
$ws.Range("D6").Value = 1.013759598832091
$ws.Range("E6").Value = 1.013261212756077
$ws.Range("F6").Value = 1.010785438287014
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.015825496594373
$ws.Range("K6").Value = 1.016239186605169
$ws.Range("L6").Value = 1.01574208822494
$ws.Range("M6").Value = 1.013272730806827
$ws.Range("N6").Value = 1.008940860364173
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.010463766062008
$ws.Range("D7").Value = 1.012962389824796
$ws.Range("E7").Value = 1.012541032099205
$ws.Range("F7").Value = 1.009748454240761
$ws.Range("I7").Value = 1.02344298551034
$ws.Range("J7").Value = 1.015119132104624
$ws.Range("K7").Value = 1.01550533668602
$ws.Range("L7").Value = 1.015085097703793
$ws.Range("M7").Value = 1.01229995991101
$ws.Range("N7").Value = 1.008689026894543
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.006991305242733
$ws.Range("D8").Value = 1.009606892739674
$ws.Range("E8").Value = 1.009508635157863
$ws.Range("F8").Value = 1.00538585128034
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.012142871873047
$ws.Range("K8").Value = 1.012414259729833
$ws.Range("L8").Value = 1.01231629258539
$ws.Range("M8").Value = 1.008205750865691
$ws.Range("N8").Value = 1.007624777058743
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.000774926865738
$ws.Range("D9").Value = 1.003601057654389
$ws.Range("E9").Value = 1.004076946789063
$ws.Range("F9").Value = 0.9975848337851936
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.006803877021899
$ws.Range("K9").Value = 1.006873078252644
$ws.Range("L9").Value = 1.007347303589047
$ws.Range("M9").Value = 1.000878039125616
$ws.Range("N9").Value = 1.005703893909823
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 0.996559851789255
$ws.Range("D10").Value = 0.9995295930248139
$ws.Range("E10").Value = 1.000391984513435
$ws.Range("F10").Value = 0.9923007265657129
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.003176380503838
$ws.Range("K10").Value = 1.003110765348435
$ws.Range("L10").Value = 1.003969824468226
$ws.Range("M10").Value = 0.9959100850204839
$ws.Range("N10").Value = 1.004391033268429
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 0.9947166684116076
$ws.Range("D11").Value = 0.9977494253751147
$ws.Range("E11").Value = 0.998780179401472
$ws.Range("F11").Value = 0.9899912410929461
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.001588398655512
$ws.Range("K11").Value = 1.0014643740569
$ws.Range("L11").Value = 1.002490978865179
$ws.Range("M11").Value = 0.9937377287925274
$ws.Range("N11").Value = 1.003814516701121
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 0.9940292131350515
$ws.Range("D12").Value = 0.9970855071619025
$ws.Range("E12").Value = 0.9981789596307075
$ws.Range("F12").Value = 0.9891300330018727
$ws.Range("I12").Value = 1.023938164268118
$ws.Range("J12").Value = 1.000995865763106
$ws.Range("K12").Value = 1.00085013801553
$ws.Range("L12").Value = 1.001939123010979
$ws.Range("M12").Value = 0.9929274990880808
$ws.Range("N12").Value = 1.003599131600842
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 0.9941768037503439
$ws.Range("D13").Value = 0.997228043017738
$ws.Range("E13").Value = 0.9983080390005967
$ws.Range("F13").Value = 0.98931491955115
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.001123089153442
$ws.Range("K13").Value = 1.000982017172391
$ws.Range("L13").Value = 1.002057614631269
$ws.Range("M13").Value = 0.9931014485550279
$ws.Range("N13").Value = 1.003645389149349
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 0.9946599011945074
$ws.Range("D14").Value = 0.9976946010513843
$ws.Range("E14").Value = 0.998730534404076
$ws.Range("F14").Value = 0.9899201227947861
$ws.Range("I14").Value = 1.023922712353277
$ws.Range("J14").Value = 1.001539475043137
$ws.Range("K14").Value = 1.001413656617625
$ws.Range("L14").Value = 1.002445414751849
$ws.Range("M14").Value = 0.9936708234728923
$ws.Range("N14").Value = 1.003796738418048
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 0.9949571772351194
$ws.Range("D15").Value = 0.9979817040314944
$ws.Range("E15").Value = 0.99899051060083
$ws.Range("F15").Value = 0.990292558607224
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.001795665392204
$ws.Range("K15").Value = 1.001679244150522
$ws.Range("L15").Value = 1.002684011100385
$ws.Range("M15").Value = 0.994021190128586
$ws.Range("N15").Value = 1.00388982421385
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 0.9966817889603912
$ws.Range("D16").Value = 0.9996473660161398
$ws.Range("E16").Value = 1.000498605863377
$ws.Range("F16").Value = 0.9924535359026343
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.003281398206676
$ws.Range("K16").Value = 1.003219658601091
$ws.Range("L16").Value = 1.004067618246909
$ws.Range("M16").Value = 0.9960537990814294
$ws.Range("N16").Value = 1.004429122571201
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 0.9977586915948756
$ws.Range("D17").Value = 1.000687517331303
$ws.Range("E17").Value = 1.00144019575883
$ws.Range("F17").Value = 0.9938032226934082
$ws.Range("I17").Value = 1.023843013862359
$ws.Range("J17").Value = 1.004208674963359
$ws.Range("K17").Value = 1.004181225239356
$ws.Range("L17").Value = 1.0049310741052
$ws.Range("M17").Value = 0.9973230313618544
$ws.Range("N17").Value = 1.004765234606172
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 0.9983850985791206
$ws.Range("D18").Value = 1.001292567750501
$ws.Range("E18").Value = 1.00198785318273
$ws.Range("F18").Value = 0.9945884141047752
$ws.Range("I18").Value = 1.023826118601775
$ws.Range("J18").Value = 1.004747881723711
$ws.Range("K18").Value = 1.004740429677374
$ws.Range("L18").Value = 1.005433139341395
$ws.Range("M18").Value = 0.9980613169681093
$ws.Range("N18").Value = 1.004960509703679
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 0.9985983965555326
$ws.Range("D19").Value = 1.001498597077148
$ws.Range("E19").Value = 1.002174328996299
$ws.Range("F19").Value = 0.9948557992536898
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.004931458856658
$ws.Range("K19").Value = 1.004930825056832
$ws.Range("L19").Value = 1.005604066228441
$ws.Range("M19").Value = 0.9983127119490636
$ws.Range("N19").Value = 1.005026963285528
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 0.9976433300527142
$ws.Range("D20").Value = 1.000576090560033
$ws.Range("E20").Value = 1.001339333670627
$ws.Range("F20").Value = 0.9936586279646489
$ws.Range("I20").Value = 1.023846096030143
$ws.Range("J20").Value = 1.004109359115976
$ws.Range("K20").Value = 1.004078230742145
$ws.Range("L20").Value = 1.004838596883239
$ws.Range("M20").Value = 0.9971870663168685
$ws.Range("N20").Value = 1.004729253196775
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 0.9945177196130279
$ws.Range("D21").Value = 0.9975572862761278
$ws.Range("E21").Value = 0.9986061904447094
$ws.Range("F21").Value = 0.9897419995063138
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.001416934683181
$ws.Range("K21").Value = 1.001286624685648
$ws.Range("L21").Value = 1.002331288293
$ws.Range("M21").Value = 0.9935032494993614
$ws.Range("N21").Value = 1.003752204366217
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 0.9925361831080397
$ws.Range("D22").Value = 0.9956436591728159
$ws.Range("E22").Value = 0.996873107129798
$ws.Range("F22").Value = 0.9872599324566306
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 0.9997085132618043
$ws.Range("K22").Value = 0.9995158005207392
$ws.Range("L22").Value = 1.000740064360185
$ws.Range("M22").Value = 0.9911678100289011
$ws.Range("N22").Value = 1.003130696647748
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 0.9935882184135105
$ws.Range("D23").Value = 0.9966596209353737
$ws.Range("E23").Value = 0.9977932668817246
$ws.Range("F23").Value = 0.9885776228259355
$ws.Range("I23").Value = 1.023948818288664
$ws.Range("J23").Value = 1.000615689345496
$ws.Range("K23").Value = 1.00045606237576
$ws.Range("L23").Value = 1.001585032801825
$ws.Range("M23").Value = 0.9924077441135375
$ws.Range("N23").Value = 1.003460863377317
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 0.9976954623140015
$ws.Range("D24").Value = 1.000626444621296
$ws.Range("E24").Value = 1.001384913709678
$ws.Range("F24").Value = 0.9937239704495554
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.004154240790481
$ws.Range("K24").Value = 1.004124774649397
$ws.Range("L24").Value = 1.004880388217769
$ws.Range("M24").Value = 0.9972485093517669
$ws.Range("N24").Value = 1.004745514034718
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.002394095298739
$ws.Range("D25").Value = 1.005165242897434
$ws.Range("E25").Value = 1.005492080668991
$ws.Range("F25").Value = 0.9996157571292119
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.008195799769047
$ws.Range("K25").Value = 1.00831726685016
$ws.Range("L25").Value = 1.008643008013573
$ws.Range("M25").Value = 1.002786521031373
$ws.Range("N25").Value = 1.006206052205764
